# Visual Elements System Test Plan Update
# PBI-009 As a player, I want to be able to see hearts based on how many
# lives I have left so that I know when I am close to game over.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Heading: collapse the three split runs ("Visual Element
#    Verification" + " " + "- Buc Battle Game") into a single run by
#    re-finding the full heading text and replacing it with itself -
#    Word's Find/Replace coalesces the matched range into one run.
# ---------------------------------------------------------------------
$dash = [char]0x2013
$headingText = "Visual Element Verification " + $dash + " Buc Battle Game"
$d.Content.Find.Execute($headingText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $headingText, 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Add the new V-3 row to the "V - Verifications/Actions" table
#    (the 2nd table in the document).
# ---------------------------------------------------------------------
$verifTable = $d.Tables(2)
$newRow = $verifTable.Rows.Add()
$newRow.Cells(1).Range.Text = "V-3"
$newRow.Cells(2).Range.Text = "Verify that hearts are displayed to represent the number of lives the player has remaining"

# ---------------------------------------------------------------------
# 3) Merge the two runs of the "- V-2: ..." expected-results paragraph
#    into a single run, the same way as step 1.
# ---------------------------------------------------------------------
$v2Text = "- V-2: One or more enemy images with pirate ship themes appear at the top of the screen."
$d.Content.Find.Execute($v2Text, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $v2Text, 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Turn the empty paragraph right after the "Actual Results" table
#    into a paragraph with three runs: " ", " ", "V-3".
#    Locate it via the "Actual Results" table's trailing empty
#    paragraph, then rebuild its content with raw WordOpenXML so the
#    three runs stay distinct (a plain Range.Text assignment would
#    collapse them back into one run).
# ---------------------------------------------------------------------
$actualResultsTable = $d.Tables(3)
$afterTableRange = $d.Range($actualResultsTable.Range.End, $actualResultsTable.Range.End)
$targetPara = $afterTableRange.Paragraphs(1)

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
'<pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
'<w:body>' + `
'<w:p>' + `
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
'<w:r><w:t>V-3</w:t></w:r>' + `
'</w:p>' + `
'</w:body>' + `
'</w:document>' + `
'</pkg:xmlData></pkg:part></pkg:package>'

$targetPara.Range.InsertXML($packageXml) | Out-Null

# ---------------------------------------------------------------------
# 5) Merge the lone page-break paragraph with the following "Sign-Off"
#    paragraph into a single paragraph (deleting the paragraph mark
#    between them keeps both runs intact under one <w:p>).
# ---------------------------------------------------------------------
$pageBreakRange = $d.Content.Find.Execute("Sign-Off") | Out-Null
$signOffFind = $d.Content
$signOffFind.Find.Execute("Sign-Off") | Out-Null
$signOffPara = $signOffFind.Paragraphs(1)
$prevPara = $signOffPara.Previous()
$markRange = $d.Range($prevPara.Range.End - 1, $prevPara.Range.End)
$markRange.Delete()
